# Update gh-pages output data (generated at 456a3b4)
# Applies numeric updates (want-to-go counts / min price) across the
# 展览 / 演出 / 本地生活 / 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G4").Value  = 168
$ws1.Range("F7").Value  = 379
$ws1.Range("F8").Value  = 69
$ws1.Range("F10").Value = 26
$ws1.Range("F11").Value = 662
$ws1.Range("F12").Value = 1507
$ws1.Range("F13").Value = 5859
$ws1.Range("F14").Value = 96
$ws1.Range("F15").Value = 1647
$ws1.Range("F16").Value = 404
$ws1.Range("F17").Value = 5561
$ws1.Range("F18").Value = 100
$ws1.Range("F22").Value = 1578
$ws1.Range("F23").Value = 827
$ws1.Range("F24").Value = 32
$ws1.Range("F25").Value = 95
$ws1.Range("F26").Value = 1164
$ws1.Range("F28").Value = 158
$ws1.Range("F29").Value = 14
$ws1.Range("F31").Value = 3832

# --- Sheet 2: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 104
$ws2.Range("F5").Value = 200
$ws2.Range("F8").Value = 313

# --- Sheet 3: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9453
$ws3.Range("F4").Value = 2180
$ws3.Range("F5").Value = 545

# --- Sheet 4: 全部类型 (All Types, merged view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9453
$ws4.Range("F4").Value  = 2180
$ws4.Range("G6").Value  = 168
$ws4.Range("F7").Value  = 545
$ws4.Range("F9").Value  = 379
$ws4.Range("F10").Value = 69
$ws4.Range("F14").Value = 662
$ws4.Range("F15").Value = 1507
$ws4.Range("F16").Value = 5859
$ws4.Range("F17").Value = 96
$ws4.Range("F18").Value = 313
$ws4.Range("F19").Value = 1647
$ws4.Range("F22").Value = 404
$ws4.Range("F25").Value = 5561
$ws4.Range("F26").Value = 100
$ws4.Range("F30").Value = 1578
$ws4.Range("F31").Value = 827
$ws4.Range("F32").Value = 32
$ws4.Range("F33").Value = 95
$ws4.Range("F34").Value = 1164
$ws4.Range("F36").Value = 158
$ws4.Range("F40").Value = 14
$ws4.Range("F45").Value = 3832
